$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Date: 2024-04-18T13:54:15+00:00 -> 2024-06-11T08:08:31+00:00
$meta.Range("B8").Value = "2024-06-11T08:08:31+00:00"

# Description: update wording to French-capitalised terms
$meta.Range("B11").Value = "CodeSystem for french prescription category (Création, Arrêt, Modification, Validation)."

# --- Concepts sheet updates ---
$concepts = $wb.Worksheets.Item("Concepts")

# Row for code "C": Display Creation -> Création
$concepts.Range("C2").Value = "Création"

# Row for code "S"/"Stop" -> code "A"/"Arrêt"
$concepts.Range("B3").Value = "A"
$concepts.Range("C3").Value = "Arrêt"
